{"js": "// Update the date label and all two-digit multiplication problems in the table.\nconst replacements = [\n  [\"2023-10-24 Tuesday\", \"2023-10-25 Wednesday\"],\n  [\"75\u00d794=\", \"13\u00d774=\"],\n  [\"57\u00d798=\", \"69\u00d799=\"],\n  [\"47\u00d731=\", \"82\u00d744=\"],\n  [\"29\u00d793=\", \"97\u00d771=\"],\n  [\"77\u00d732=\", \"90\u00d765=\"],\n  [\"80\u00d719=\", \"83\u00d752=\"],\n  [\"59\u00d726=\", \"72\u00d714=\"],\n  [\"89\u00d771=\", \"74\u00d764=\"],\n  [\"19\u00d751=\", \"63\u00d769=\"],\n  [\"82\u00d743=\", \"88\u00d791=\"],\n  [\"20\u00d721=\", \"74\u00d759=\"],\n  [\"16\u00d772=\", \"31\u00d767=\"],\n  [\"83\u00d784=\", \"66\u00d785=\"],\n  [\"90\u00d776=\", \"11\u00d771=\"],\n  [\"72\u00d796=\", \"55\u00d753=\"],\n  [\"26\u00d790=\", \"35\u00d797=\"],\n  [\"91\u00d731=\", \"86\u00d721=\"],\n  [\"62\u00d745=\", \"44\u00d724=\"],\n  [\"33\u00d737=\", \"15\u00d793=\"],\n  [\"11\u00d792=\", \"98\u00d777=\"],\n  [\"26\u00d756=\", \"72\u00d714=\"],\n  [\"69\u00d742=\", \"76\u00d777=\"],\n  [\"19\u00d723=\", \"71\u00d741=\"],\n  [\"80\u00d756=\", \"95\u00d792=\"],\n  [\"11\u00d733=\", \"24\u00d798=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace only the first occurrence, since every \"before\" string in this\n  // document is unique.\n  found.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date label and all two-digit multiplication problems in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-10-24 Tuesday\", \"2023-10-25 Wednesday\"),\n    @(\"75\u00d794=\", \"13\u00d774=\"),\n    @(\"57\u00d798=\", \"69\u00d799=\"),\n    @(\"47\u00d731=\", \"82\u00d744=\"),\n    @(\"29\u00d793=\", \"97\u00d771=\"),\n    @(\"77\u00d732=\", \"90\u00d765=\"),\n    @(\"80\u00d719=\", \"83\u00d752=\"),\n    @(\"59\u00d726=\", \"72\u00d714=\"),\n    @(\"89\u00d771=\", \"74\u00d764=\"),\n    @(\"19\u00d751=\", \"63\u00d769=\"),\n    @(\"82\u00d743=\", \"88\u00d791=\"),\n    @(\"20\u00d721=\", \"74\u00d759=\"),\n    @(\"16\u00d772=\", \"31\u00d767=\"),\n    @(\"83\u00d784=\", \"66\u00d785=\"),\n    @(\"90\u00d776=\", \"11\u00d771=\"),\n    @(\"72\u00d796=\", \"55\u00d753=\"),\n    @(\"26\u00d790=\", \"35\u00d797=\"),\n    @(\"91\u00d731=\", \"86\u00d721=\"),\n    @(\"62\u00d745=\", \"44\u00d724=\"),\n    @(\"33\u00d737=\", \"15\u00d793=\"),\n    @(\"11\u00d792=\", \"98\u00d777=\"),\n    @(\"26\u00d756=\", \"72\u00d714=\"),\n    @(\"69\u00d742=\", \"76\u00d777=\"),\n    @(\"19\u00d723=\", \"71\u00d741=\"),\n    @(\"80\u00d756=\", \"95\u00d792=\"),\n    @(\"11\u00d733=\", \"24\u00d798=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
